$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextCell "D2" "62.988.94"
Set-TextCell "E2" "  +0.01%  "
Set-TextCell "D3" "2.543.28"
Set-TextCell "E3" "  +3.67%  "
Set-TextCell "E4" "  -0.04%  "
Set-TextCell "D5" "568.80"
Set-TextCell "E5" "  +1.00%  "
Set-TextCell "D6" "145.89"
Set-TextCell "E6" "  +2.76%  "
Set-TextCell "E7" "  -0.02%  "
Set-TextCell "D8" "0.583"
Set-TextCell "E8" "  +0.24%  "
Set-TextCell "D9" "2.541.19"
Set-TextCell "E9" "  +3.66%  "
Set-TextCell "D10" "0.105"
Set-TextCell "E10" "  +0.26%  "
Set-TextCell "D11" "5.50"
Set-TextCell "E11" "  -2.72%  "
Set-TextCell "E12" "  -0.03%  "
Set-TextCell "D13" "0.352"
Set-TextCell "E13" "  -0.13%  "
Set-TextCell "D14" "27.26"
Set-TextCell "E14" "  +1.25%  "
Set-TextCell "D15" "2.996.24"
Set-TextCell "E15" "  +3.61%  "
Set-TextCell "D16" "62.887.23"
Set-TextCell "E16" "  +0.04%  "
Set-TextCell "E17" "  +1.30%  "
Set-TextCell "D18" "2.548.00"
Set-TextCell "E18" "  +3.92%  "
Set-TextCell "D19" "11.32"
Set-TextCell "E19" "  +0.97%  "
Set-TextCell "E20" "  +1.67%  "
Set-TextCell "D21" "333.50"
Set-TextCell "E21" "  -1.56%  "
Set-TextCell "D22" "6.81"
Set-TextCell "E22" "  +1.00%  "
Set-TextCell "D23" "0.999"
Set-TextCell "E23" "  -0.05%  "
Set-TextCell "D24" "65.07"
Set-TextCell "E24" "  -0.24%  "
Set-TextCell "E25" "  -0.74%  "
Set-TextCell "D26" "1.60"
Set-TextCell "E26" "  +7.69%  "
Set-TextCell "E27" "  -0.04%  "
Set-TextCell "B28" "SuiNetwork"
Set-TextCell "C28" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell "D28" "1.47"
Set-TextCell "E28" "  +3.24%  "
Set-TextCell "B29" "InternetComputer(DFINITY)"
Set-TextCell "C29" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D29" "8.36"
Set-TextCell "E29" "  +4.38%  "
Set-TextCell "D30" "7.25"
Set-TextCell "E30" "  +7.85%  "
Set-TextCell "D31" "0.0₃0813"
Set-TextCell "E31" "  +3.59%  "
Set-TextCell "D32" "1.85"
Set-TextCell "E32" "  +0.84%  "
Set-TextCell "D33" "175.35"
Set-TextCell "E33" "  -0.81%  "
Set-TextCell "E34" "  +3.27%  "
Set-TextCell "D35" "409.21"
Set-TextCell "E35" "  +6.84%  "
Set-TextCell "D36" "0.399"
Set-TextCell "E36" "  +0.54%  "
Set-TextCell "D37" "18.95"
Set-TextCell "E37" "  +1.34%  "
Set-TextCell "E38" "  -0.01%  "
Set-TextCell "D39" "4.35"
Set-TextCell "E39" "  +0.81%  "
Set-TextCell "D40" "1.75"
Set-TextCell "E40" "  +2.11%  "
Set-TextCell "E41" "  -0.03%  "
Set-TextCell "D42" "39.66"
Set-TextCell "E42" "  -0.73%  "
Set-TextCell "D43" "151.98"
Set-TextCell "E43" "  +1.98%  "
Set-TextCell "E44" "  +1.79%  "
Set-TextCell "D45" "20.75"
Set-TextCell "E45" "  +1.72%  "
Set-TextCell "D46" "0.602"
Set-TextCell "E46" "  +1.35%  "
Set-TextCell "D47" "0.0529"
Set-TextCell "E47" "  +2.98%  "
Set-TextCell "E48" "  +0.28%  "
Set-TextCell "D49" "0.0238"
Set-TextCell "E49" "  +4.17%  "
Set-TextCell "E50" "  +2.26%  "
Set-TextCell "E51" "  -1.23%  "
